$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Shift the four existing "role" rows (old rows 4-7) down by one row to make
# room for the new HEALED-consortium entry, preserving each destination
# row's own height (rows keep their position-based height, exactly like
# manually re-typing shifted data into place) while carrying per-cell
# formatting down with the content. Must proceed bottom-to-top so we never
# clobber a source row before it has been copied down.
# ---------------------------------------------------------------------------

$ws.Range("A7:E7").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A6:E6").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)

$ws.Range("A5:E5").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)

$ws.Range("A4:E4").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# New row 4's E cell uses the wrap-only style (same as the header/"why"
# column cells in rows 1-3), not the style row 4 used to carry.
$ws.Range("E1").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Row 8 (now holding the old row 7 "Post-doctoral Researcher" entry) keeps
# its own pre-existing height (70.5) automatically since PasteSpecial only
# moves formats, never row heights. Rows 5-7 likewise keep the heights that
# already lived at those row numbers. Row 9 is untouched.
# ---------------------------------------------------------------------------

# Shift values (rows 5-8 now hold what rows 4-7 used to hold).
$ws.Cells.Item(8,1).Value2 = $ws.Cells.Item(7,1).Value2
$ws.Cells.Item(8,2).Value2 = $ws.Cells.Item(7,2).Value2
$ws.Cells.Item(8,3).Value2 = $ws.Cells.Item(7,3).Value2
$ws.Cells.Item(8,4).Value2 = $ws.Cells.Item(7,4).Value2
$ws.Cells.Item(8,5).Value2 = $ws.Cells.Item(7,5).Value2

$ws.Cells.Item(7,1).Value2 = $ws.Cells.Item(6,1).Value2
$ws.Cells.Item(7,2).Value2 = $ws.Cells.Item(6,2).Value2
$ws.Cells.Item(7,3).Value2 = $ws.Cells.Item(6,3).Value2
$ws.Cells.Item(7,4).Value2 = $ws.Cells.Item(6,4).Value2
$ws.Cells.Item(7,5).Value2 = $ws.Cells.Item(6,5).Value2

$ws.Cells.Item(6,1).Value2 = $ws.Cells.Item(5,1).Value2
$ws.Cells.Item(6,2).Value2 = $ws.Cells.Item(5,2).Value2
$ws.Cells.Item(6,3).Value2 = $ws.Cells.Item(5,3).Value2
$ws.Cells.Item(6,4).Value2 = $ws.Cells.Item(5,4).Value2
$ws.Cells.Item(6,5).Value2 = $ws.Cells.Item(5,5).Value2

$ws.Cells.Item(5,1).Value2 = $ws.Cells.Item(4,1).Value2
$ws.Cells.Item(5,2).Value2 = $ws.Cells.Item(4,2).Value2
$ws.Cells.Item(5,3).Value2 = $ws.Cells.Item(4,3).Value2
$ws.Cells.Item(5,4).Value2 = $ws.Cells.Item(4,4).Value2
$ws.Cells.Item(5,5).Value2 = $ws.Cells.Item(4,5).Value2

# New row 4: the HEALED consortium entry (reuses "Research Fellow" / "Trinity
# College Dublin" text already used elsewhere on the sheet).
$ws.Cells.Item(4,1).Value2 = "Research Fellow"
$ws.Cells.Item(4,2).Value2 = "Supervisor: Aideen Long"
$ws.Cells.Item(4,3).Value2 = "2023-2024"
$ws.Cells.Item(4,4).Value2 = "Trinity College Dublin"
$ws.Cells.Item(4,5).Value2 = "HEALED consortium - molecular biology and sequencing lead of tumour DNA and RNA-sequencing (with Aideen Long), since August 2023"

# ---------------------------------------------------------------------------
# Sheet view: the live selection moved to E4 and the sheet no longer pins a
# frozen/scrolled top-left cell (was A4).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("E4").Select()
